$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25, pushing existing rows 25-32 down to 26-33.
$ws.Rows("25:25").Insert()

# Populate the newly inserted row 25 with data (same record as the old row 25,
# but with an updated date and volume, matching a new weekly price entry).
$ws.Cells.Item(25, 1).Value = 2
$ws.Cells.Item(25, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(25, 3).Value = "Coquimbo"
$ws.Cells.Item(25, 4).Value = 44455
$ws.Cells.Item(25, 5).Value = 4
$ws.Cells.Item(25, 6).Value = 100112022
$ws.Cells.Item(25, 7).Value = "Arveja Verde"
$ws.Cells.Item(25, 8).Value = "Perfection"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 800
$ws.Cells.Item(25, 11).Value = 28000
$ws.Cells.Item(25, 12).Value = 30000
$ws.Cells.Item(25, 13).Value = 29000
$ws.Cells.Item(25, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(25, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(25, 16).Value = 1160
$ws.Cells.Item(25, 17).Value = 25
$ws.Cells.Item(25, 18).Value = "Hortaliza"
